$wb = $excel.ActiveWorkbook

# Ordered list of new sheet names matching the existing sheet order (1-based index)
$newNames = @(
    "summ35338388",
    "summ35563942",
    "summ35824308",
    "summ36075145",
    "summ36316770",
    "summ36519946",
    "summ36717514",
    "summ36920350",
    "summ37124477",
    "summ37324479",
    "summ37519398",
    "summ37714847",
    "summ37925511",
    "summ38137079",
    "summ38389275",
    "summ38611238",
    "summ38813865",
    "summ39015734",
    "summ39265280",
    "summ39459565",
    "summ39655417",
    "summ39862186",
    "summ40056728",
    "summ40261756",
    "summ40462457",
    "summ40664903",
    "summ40870602",
    "summ41073953",
    "summ41277971",
    "summ41477182",
    "summ41673552",
    "summ41870424",
    "summ42074638",
    "summ42269508",
    "summ42469808",
    "summ42670374",
    "summ42968835",
    "summ43189391",
    "summ43391898",
    "summ43592520",
    "summ43794620",
    "summ43988753",
    "summ44199151",
    "summ44435607",
    "summ44704796",
    "summ45016578",
    "summ45289204",
    "summ45553748",
    "summ45817306",
    "summ46105717"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

